# Natmi following Dr Hou advice
# Rebuild the LR-pairs data table to include the full 3x3 cross product of
# clusters (ECs, FAPs, sCs) as Sending cluster x Target cluster, with
# updated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs",  "Hgf", "Sdc2", "ECs",  2, 0.6666666666666666, 5.928568666666667, 17.785706,          0.3809768389628236, 0.3809768389628236, 3, 1, 0.7065936666666666, 2.119781,           0.005187843618793344, 0.005187843618793344, 4.18908907226511,   37.701801650386,    0.001976448262921344, 0.001976448262921344),
  @("ECs",  "Hgf", "Sdc2", "FAPs", 2, 0.6666666666666666, 5.928568666666667, 17.785706,          0.3809768389628236, 0.3809768389628236, 3, 1, 111.9320066666667,  335.79602,           0.8218100075305903,   0.8218100075305903,   663.5965875211244,  5972.36928769012,   0.3130905788970186,   0.3130905788970186),
  @("ECs",  "Hgf", "Sdc2", "sCs",  2, 0.6666666666666666, 5.928568666666667, 17.785706,          0.3809768389628236, 0.3809768389628236, 3, 1, 23.563205,          70.689615,           0.1730021488506163,   0.1730021488506163,   139.6960788492433,  1257.26470964319,   0.06590981180288369,  0.0659098118028837),
  @("FAPs", "Hgf", "Sdc2", "ECs",  3, 1,                  8.583520999999999, 25.750563,          0.5515872180307627, 0.5515872180307626, 3, 1, 0.7065936666666666, 2.119781,           0.005187843618793344, 0.005187843618793344, 6.065061576300332,  54.58555418670299,  0.002861548229268865, 0.002861548229268865),
  @("FAPs", "Hgf", "Sdc2", "FAPs", 3, 1,                  8.583520999999999, 25.750563,          0.5515872180307627, 0.5515872180307626, 3, 1, 111.9320066666667,  335.79602,           0.8218100075305903,   0.8218100075305903,   960.7707297954732,  8646.936568159259,  0.4532998958036384,   0.4532998958036384),
  @("FAPs", "Hgf", "Sdc2", "sCs",  3, 1,                  8.583520999999999, 25.750563,          0.5515872180307627, 0.5515872180307626, 3, 1, 23.563205,          70.689615,           0.1730021488506163,   0.1730021488506163,   202.255264944805,   1820.297384503245,  0.09542577399785536,  0.09542577399785535),
  @("sCs",  "Hgf", "Sdc2", "ECs",  3, 1,                  1.049404,           3.148212,           0.06743594300641363, 0.06743594300641362, 3, 1, 0.7065936666666666, 2.119781,           0.005187843618793344, 0.005187843618793344, 0.7415022201746666, 6.673519981571999,  0.0003498471266031346,0.0003498471266031345),
  @("sCs",  "Hgf", "Sdc2", "FAPs", 3, 1,                  1.049404,           3.148212,           0.06743594300641363, 0.06743594300641362, 3, 1, 111.9320066666667,  335.79602,           0.8218100075305903,   0.8218100075305903,   117.4618955240267,  1057.15705971624,   0.05541953282993325,  0.05541953282993323),
  @("sCs",  "Hgf", "Sdc2", "sCs",  3, 1,                  1.049404,           3.148212,           0.06743594300641363, 0.06743594300641362, 3, 1, 23.563205,          70.689615,           0.1730021488506163,   0.1730021488506163,   24.72732157982,     222.54589421838,    0.01166656304987725,  0.01166656304987725)
)

$rowIndex = 2
foreach ($row in $data) {
  $colIndex = 1
  foreach ($val in $row) {
    $ws.Cells.Item($rowIndex, $colIndex).Value = $val
    $colIndex++
  }
  $rowIndex++
}
